$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.976.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.18%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.419.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.56%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'578.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.65%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'152.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.39%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.06%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +1.24%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +4.03%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -0.81%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +3.29%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.002.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.58%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.77%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'28.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.63%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.432.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.21%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -0.45%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'62.021.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.03%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +1.66%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'14.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.63%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -4.42%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'382.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.75%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.572"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.01%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'75.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.87%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.10%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'3.560.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.53%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -3.65%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -1.29%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.19%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.18%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'7.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.73%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.16%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.12%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -1.42%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'23.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.98%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +3.17%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.49%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'6.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.94%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'168.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.30%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'30.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.66%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'3.452.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.68%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +2.40%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'42.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.53%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.781"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.34%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'4.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.75%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -3.66%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D47").Value = "'2.552.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.72%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'6.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.58%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'22.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.96%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -5.77%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.06%  "
$ws.Range("E51").Style = "Normal"
